# Bugfixed the naive forecaster component module
#
# Column A held raw date serials (formatted as YYYY-MM-DD HH:MM:SS) for each
# observation. Replace them with plain "YYYYQn" quarter-label text, matching
# the header's style (center/top aligned, bold, bordered, General format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerCell = $ws.Range("A1")
$epoch = Get-Date -Year 1899 -Month 12 -Day 30

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = $ws.UsedRange.Rows.Count }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = [double]$cell.Value2
    $dt = $epoch.AddDays($serial)
    $quarter = [Math]::Floor(($dt.Month - 1) / 3) + 1
    $label = "{0}Q{1}" -f $dt.Year, $quarter

    $cell.Value = $label

    $headerCell.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = 0
